# ---------------------------------------------------------------------------
# Sheaft catalogue produits - workbook update
#  - Simplify the "Catalogue" table headers (drop the parenthetical hints)
#  - Remove the "Poids total" column
#  - Add a new hidden "Listes" sheet holding the reference lists (as Excel
#    Tables) used to drive dropdown data-validation on the Catalogue sheet
#  - Add workbook-scoped defined names pointing at those list tables
#  - Wire data validation (lists + numeric) onto the Catalogue columns
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Rename the Catalogue headers (text simplified, hints removed) -----
$ws.Range("D1").Value = "TVA % "
$ws.Range("E1").Value = "Conditionnement"
$ws.Range("G1").Value = "Unité de mesure`nUniquement si conditionnement = Poids"
$ws.Range("H1").Value = "Catégorie"
$ws.Range("I1").Value = "Agriculture Bio"
$ws.Range("J1").Value = "Disponible à la vente "

# --- 2. Drop the trailing "Poids total" column -----------------------------
$lo.ListColumns.Item(12).Delete()
$ws.Range("L1").Clear()

# --- 3. Add the hidden "Listes" sheet with the reference lists ------------
$wsListes = $wb.Worksheets.Add($null, $ws)
$wsListes.Name = "Listes"

# TVA (A1:A4)
$wsListes.Cells.Item(1,1).Value = "TVA"
$wsListes.Cells.Item(2,1).Value = 5.5
$wsListes.Cells.Item(3,1).Value = 10
$wsListes.Cells.Item(4,1).Value = 20

# Conditionnement (C1:C6)
$wsListes.Cells.Item(1,3).Value = "Conditionnement"
$wsListes.Cells.Item(2,3).Value = "Boîte"
$wsListes.Cells.Item(3,3).Value = "Botte"
$wsListes.Cells.Item(4,3).Value = "Bouquet"
$wsListes.Cells.Item(5,3).Value = "Pièce"
$wsListes.Cells.Item(6,3).Value = "Poids"

# Unité de mesure (E1:E5)
$wsListes.Cells.Item(1,5).Value = "Unité de mesure"
$wsListes.Cells.Item(2,5).Value = "ml"
$wsListes.Cells.Item(3,5).Value = "L"
$wsListes.Cells.Item(4,5).Value = "g"
$wsListes.Cells.Item(5,5).Value = "kg"

# Catégorie (G1:G7)
$wsListes.Cells.Item(1,7).Value = "Catégorie"
$wsListes.Cells.Item(2,7).Value = "Boisson"
$wsListes.Cells.Item(3,7).Value = "Épicerie"
$wsListes.Cells.Item(4,7).Value = "Fruits et légumes"
$wsListes.Cells.Item(5,7).Value = "Oeufs et produits laitiers"
$wsListes.Cells.Item(6,7).Value = "Poisson"
$wsListes.Cells.Item(7,7).Value = "Viande"

# Boolean (I1:I3)
$wsListes.Cells.Item(1,9).Value = "Boolean"
$wsListes.Cells.Item(2,9).Value = "Oui"
$wsListes.Cells.Item(3,9).Value = "Non"

# --- 4. Turn each list range into an Excel Table ---------------------------
$loTva = $wsListes.ListObjects.Add(1, $wsListes.Range("A1:A4"), $null, 1)
$loTva.Name = "Tableau1"
$loTva.TableStyle = "TableStyleLight13"

$loCond = $wsListes.ListObjects.Add(1, $wsListes.Range("C1:C6"), $null, 1)
$loCond.Name = "Tableau2"
$loCond.TableStyle = "TableStyleLight13"

$loUm = $wsListes.ListObjects.Add(1, $wsListes.Range("E1:E5"), $null, 1)
$loUm.Name = "Tableau3"
$loUm.TableStyle = "TableStyleLight13"

$loCat = $wsListes.ListObjects.Add(1, $wsListes.Range("G1:G7"), $null, 1)
$loCat.Name = "Tableau5"
$loCat.TableStyle = "TableStyleLight13"

$loBool = $wsListes.ListObjects.Add(1, $wsListes.Range("I1:I3"), $null, 1)
$loBool.Name = "Tableau6"
$loBool.TableStyle = "TableStyleLight13"

# --- 5. Hide the "Listes" sheet --------------------------------------------
$wsListes.Visible = 0

# --- 6. Workbook-scoped defined names used by the data validations --------
$wb.Names.Add("liste_tva", "=Tableau1[TVA]")
$wb.Names.Add("liste_conditionnement", "=Tableau2[Conditionnement]")
$wb.Names.Add("liste_um", "=Tableau3[Unité de mesure]")
$wb.Names.Add("liste_categorie", "=Tableau5[Catégorie]")
$wb.Names.Add("liste_bool", "=Tableau6[Boolean]")

# --- 7. Data validation on the Catalogue sheet ------------------------------
$msgPrompt = "Cliquez sur la flèche à droite de la case"
$errTitle = "Erreur"
$errMsg = "Attention, il faut saisir un nombre (entier ou décimal) supérieur à 0. Cliquez sur ""Rééssayer"" pour modifier votre saisie."

$rngC = $ws.Range("C2:C1048576")
$rngC.Validation.Add(2, 1, 5, 0)
$rngC.Validation.ErrorTitle = $errTitle
$rngC.Validation.ErrorMessage = $errMsg
$rngC.Validation.ShowError = $true

$rngD = $ws.Range("D2:D1048576")
$rngD.Validation.Add(3, 1, 1, "=liste_tva")
$rngD.Validation.InputMessage = $msgPrompt
$rngD.Validation.ShowInput = $true

$rngE = $ws.Range("E2:E1048576")
$rngE.Validation.Add(3, 1, 1, "=liste_conditionnement")
$rngE.Validation.InputMessage = $msgPrompt
$rngE.Validation.ShowInput = $true

$rngF = $ws.Range("F2:F1048576")
$rngF.Validation.Add(2, 1, 5, 0)
$rngF.Validation.ErrorTitle = $errTitle
$rngF.Validation.ErrorMessage = $errMsg
$rngF.Validation.ShowError = $true

$rngG = $ws.Range("G2:G1048576")
$rngG.Validation.Add(3, 1, 1, "=liste_um")
$rngG.Validation.InputMessage = $msgPrompt
$rngG.Validation.ShowInput = $true

$rngH = $ws.Range("H2:H1048576")
$rngH.Validation.Add(3, 1, 1, "=liste_categorie")
$rngH.Validation.InputMessage = $msgPrompt
$rngH.Validation.ShowInput = $true

$rngI = $ws.Range("I2:I1048576")
$rngI.Validation.Add(3, 1, 1, "=liste_bool")
$rngI.Validation.InputMessage = $msgPrompt
$rngI.Validation.ShowInput = $true

$rngJ = $ws.Range("J2:J1048576")
$rngJ.Validation.Add(3, 1, 1, "=liste_bool")
$rngJ.Validation.InputMessage = $msgPrompt
$rngJ.Validation.ShowInput = $true

# --- 8. Misc cosmetic tweaks (column widths, selection) --------------------
$ws.Columns.Item(1).ColumnWidth = 17
$ws.Columns.Item(4).ColumnWidth = 8.1667
$ws.Columns.Item(5).ColumnWidth = 20.8333
$ws.Columns.Item(8).ColumnWidth = 22.5
$ws.Columns.Item(10).ColumnWidth = 20.6667
$ws.Columns.Item(11).ColumnWidth = 44.1667

$ws.Range("C14").Select()
